# Implemented SSN and CC Discovery
# Replace the sample source-directory rows with a single, real test path,
# drop the old UNC-share hyperlinks, and trim the sheet back down to the
# header + one data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source")

# Remove the three UNC-path hyperlinks (and their relationships) before we
# touch the rows they live on.
$ws.Hyperlinks.Delete()

# Update the remaining sample row to the new local test path...
$ws.Range("A2").Value = "C:\Users\tnabbefeld\Documents\test"

# ...and delete the now-unneeded rows 3 and 4 (adelev/meganb + meganb/admin
# sample rows), shrinking the sheet to A1:B2.
$ws.Rows("3:4").Delete()

# Match the author's final selection: the whole of row 2 selected, with A2
# as the active cell.
$ws.Rows("2:2").Select()
